$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that sits between the "${"
#    and "Property}" runs in the InputEnd row (table), merging
#    those two runs back into a single "${Property}" run. This has
#    to happen before we re-add "_GoBack" elsewhere, since this
#    runtime does not replace an existing same-named bookmark when
#    a new one is added (real Word moves it; here it would just add
#    a second one).
# -----------------------------------------------------------------
$oldGoBack = $d.Bookmarks("_GoBack")
$oldGoBackStart = $oldGoBack.Start
$bmContext = $d.Range([Math]::Max(0, $oldGoBackStart - 60), $oldGoBackStart + 60)
$bmContext.Find.Execute('${Property}', $true, $false, $false, $false, $false, $true, 1, $false, '${Property}', 1)

# -----------------------------------------------------------------
# 2) Version/date paragraph: "Version 11.08.00, 2016-01-24"
#    becomes "Version 11.08.01, 2016-02-14".
# -----------------------------------------------------------------
$verPara = $d.Paragraphs(3).Range
$verStart = $verPara.Start

# "Version 11.08.00, 2016-01-24"
#  0123456789...
# offsets are relative to the paragraph start
$patchDigit = $verStart + 15   # the "0" in ".00" -> "1"
$monthDigit = $verStart + 24   # the "1" in "-01-" -> "2"
$dayDigit   = $verStart + 26   # the "2" in "-24" -> "1"

$d.Range($patchDigit, $patchDigit + 1).Text = "1"
$d.Range($monthDigit, $monthDigit + 1).Text = "2"
$d.Range($dayDigit, $dayDigit + 1).Text = "1"

# -----------------------------------------------------------------
# 3) Word re-creates "_GoBack" at the point of the latest edit:
#    right after the "1" we just typed for the day, just before the
#    trailing "4" ("2016-02-1|4").
# -----------------------------------------------------------------
$goBackPos = $dayDigit + 1
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# -----------------------------------------------------------------
# Word keeps bookmark ids densely packed starting at 0, so removing
# the old "_GoBack" and adding the new one automatically bumps the
# pre-existing "replaceValue" bookmark from id 0 to id 1 (and its
# matching bookmarkEnd along with it) -- nothing else to do here.
# -----------------------------------------------------------------
